$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.747.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.826.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.75%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '350.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.03%  '
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.618'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0848'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.271.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.984'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.823.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.781.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +12.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.96%  '
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.13%  '
$ws.Range("E25").Value = '  +1.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.26'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '52.82'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0894'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.07%  '
$ws.Range("E35").Value = '  +0.96%  '
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.92'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.73%  '
$ws.Range("E39").Value = '  +2.05%  '
$ws.Range("E40").Value = '  +2.29%  '
$ws.Range("E41").Value = '  +1.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.90%  '
$ws.Range("E47").Value = '  +6.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.165.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.246'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +22.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.950'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.91%  '
